$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Exceptions")

# Clear the old exception table, leaving only a note pointing to the tracked doc.
$ws.Cells.Clear()
$ws.Range("A1").Value = "Tracked in FitMe/Documentation/ExceptionList.md"

[void]$ws.Range("C6").Select()
